$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)

# Explicit position/size for the content placeholder (previously inherited from layout/master).
$shp.Left = 66
$shp.Top = 143.75
$shp.Width = 828
$shp.Height = 376.4711811023622

$tr = $shp.TextFrame.TextRange
$tr.Text = ""
$lines = @(
  "Quiz Wednesday on lectures 4 and 5, excluding example interview question on summing two numbers in an array.",
  "Quizzes are not weighted; grades are averaged regardless of points used in each. E.g. 1/1, 4/8, & 0/10 give 75% after dropping lowest.",
  "See blackboard announcement on some sample projects in my GitHub account that you could start from for a portfolio project.",
  "It is acceptable to have work-in-progress projects in your portfolio.",
  "Email me if you're not sure if something belongs on GitHub.",
  "Still missing 5 GitHub usernames. Will send emails if they're not in tonight."
)
$tr.Text = [string]::Join("`r", $lines)
